$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pipe-delimited data for rows 2-24 (A..I); "<EMPTY>" marks a blank cell
$data = @"
Trivento|58448688|Cat|Cat Adoption Condo Rooms|05/06/2025|Domestic Longhair|2y 2m 1d|Hold - Behavior Mod.|62.5
CHATEAU ST. MICHELLE|58459452|Cat|Cat Adoption Condo Rooms|05/08/2025|Domestic Longhair|3y 2m 0d|Hold - Behavior Mod.|60.8
BEATRICE|58442198|Cat|Foster Home|05/06/2025|Domestic Shorthair|2y 2m 1d|In Foster|62.8
Sam|58444073|Cat|Foster Home|05/06/2025|Domestic Shorthair|3m 0d|In Foster|62.8
Masie|58450156|Cat|Foster Home|05/07/2025|Domestic Longhair|1y 8m 1d|In Foster|61.9
Quinta Nova|58454445|Cat|If The Fur Fits|05/07/2025|Domestic Longhair|2y 2m 1d|In If the Fur Fits - Behavior|61.7
Momo|58448506|Cat|If The Fur Fits|05/06/2025|Ragdoll|4m 12d|In If the Fur Fits - Medical|62.6
Elvis|58452209|Cat|If The Fur Fits|05/07/2025|Domestic Shorthair|1y 2m 1d|In If the Fur Fits - Medical|61.8
<EMPTY>|58447679|Cat|If The Fur Fits|05/06/2025|Snowshoe|3m 7d|In If the Fur Fits - Trial|62.7
Missy|58764565|Dog|Dog Adoptions B|06/23/2025|Mixed Breed, Large (over 44 lbs fully grown)|3y 14d|Hold - Meet and Greet|14.9
Karma|58764550|Dog|Dog Adoptions B|06/23/2025|Mixed Breed, Large (over 44 lbs fully grown)|2y 15d|Hold - Meet and Greet|14.9
Onyx|58622436|Dog|Dog Adoptions D|06/21/2025|Shepherd|6y 1m 6d|Available - ITFF Behavior|16.8
CINDER|58753979|Dog|Dog Holding E|06/20/2025|Mixed Breed, Large (over 44 lbs fully grown)|10y 4m 10d|Hold - Rescue|17.6
Diamond|57957792|Dog|Dog Holding F|03/04/2025|Mixed Breed, Medium (up to 44 lbs fully grown)|1y 1m 26d|Hold - Behavior|125.7
EASTON|58718199|Dog|If The Fur Fits|06/18/2025|Mixed Breed, Medium (up to 44 lbs fully grown)|8y 2m 0d|In If the Fur Fits - Medical|19.9
Magnus|58690617|Dog|If The Fur Fits|06/19/2025|Mixed Breed, Large (over 44 lbs fully grown)|4y 27d|In If the Fur Fits - Trial|18.8
JEWELIE|58753970|Dog|If The Fur Fits|06/20/2025|Mixed Breed, Small (under 24 lbs fully grown)|10y 17d|In If the Fur Fits - Trial|17.6
FIREFLY|58695418|Rabbit|Cat Adoption Condo Rooms|06/20/2025|Mix|<EMPTY>|Available|17.8
Luna|58744079|Donkey|Farm|06/17/2025|Miniature Donkey|15y 19d|Hold - Doc|20.7
Moose|58728965|Farm Type Fowl|Farm|06/17/2025|Chicken|8m 23d|Available|20.7
FUZZY SLIPPER|58757630|Rabbit|Foster Home|06/21/2025|Mini-Lop|2m 29d|In Foster|16.8
FLIP FLOP|58757640|Rabbit|Foster Home|06/21/2025|Rex|2m 29d|In Foster|16.8
Leonardo|58764920|Rabbit|Foster Home|06/23/2025|Satin|1y 15d|In Foster|14.8
"@

$lines = $data -split "`n"
$rowNum = 2
foreach ($line in $lines) {
    $cols = $line -split "\|"
    $name = $cols[0]
    $aid = $cols[1]
    $species = $cols[2]
    $location = $cols[3]
    $intakeDate = $cols[4]
    $breed = $cols[5]
    $age = $cols[6]
    $stage = $cols[7]
    $los = $cols[8]

    if ($name -eq "<EMPTY>") { $name = "" }
    if ($age -eq "<EMPTY>") { $age = "" }

    $ws.Range("A$rowNum").Value = $name
    $ws.Range("B$rowNum").Value = "'" + $aid
    $ws.Range("C$rowNum").Value = $species
    $ws.Range("D$rowNum").Value = $location
    $ws.Range("E$rowNum").Value = "'" + $intakeDate
    $ws.Range("F$rowNum").Value = $breed
    $ws.Range("G$rowNum").Value = $age
    $ws.Range("H$rowNum").Value = $stage
    $ws.Range("I$rowNum").Value = [double]$los

    $rowNum++
}
